# Re-generated project files: update PID4Cat Excel model sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "PID4CatRecord": shift/rename header columns and drop the
# trailing "change_log" column (J1), shrinking the sheet from A1:J1 to A1:I1.
$wsRecord = $wb.Worksheets.Item("PID4CatRecord")
$wsRecord.Range("D1").Value = "pid_schema_version"
$wsRecord.Range("E1").Value = "license"
$wsRecord.Range("F1").Value = "curation_contact_email"
$wsRecord.Range("G1").Value = "resource_info"
$wsRecord.Range("H1").Value = "related_identifiers"
$wsRecord.Range("I1").Value = "change_log"
$wsRecord.Range("J1").ClearContents()

# --- Sheet "ResourceInfo": extend the resource_category validation list.
$wsResourceInfo = $wb.Worksheets.Item("ResourceInfo")
$wsResourceInfo.Range("C2:C1048576").Validation.Modify(3, 1, 1, '"COLLECTION,SAMPLE,MATERIAL,DEVICE,DATA_OBJECT,DATA_SERVICE"')

# --- Sheet "LogRecord": replace RIGHTS with LICENSE in the changed_field validation list.
$wsLogRecord = $wb.Worksheets.Item("LogRecord")
$wsLogRecord.Range("C2:C1048576").Validation.Modify(3, 1, 1, '"STATUS,RESOURCE_INFO,RELATED_IDS,CONTACT,LICENSE"')

# --- Sheet "Agent": rename contact_information -> email, person_orcid -> orcid.
$wsAgent = $wb.Worksheets.Item("Agent")
$wsAgent.Range("B1").Value = "email"
$wsAgent.Range("C1").Value = "orcid"
